# Refactory like resolvers & service
# Apply the changes to the "Chat" sheet's resolver-plan table:
#  - replace the old enterChatRoom/createChatRoom/exitChatRoom rows with the
#    new getChatRoomInPost / getChatRoomInList resolvers
#  - drop the now-unused 4th data row entirely (table shrinks from A2:F6 to A2:F5)
#  - move the selection to E9 (post-edit cursor position)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chat")
$ws.Activate()

# --- Row 4: query/enterChatRoom -> mutation/getChatRoomInPost(postId) ---
$ws.Range("A4").ClearContents()
$ws.Range("B4").Value = "mutation"
$ws.Range("C4").Value = "getChatRoomInPost"
$ws.Range("E4").Value = "postId"
$ws.Range("F4").ClearContents()

# --- Row 5: mutation/createChatRoom -> query/getChatRoomInList ---
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "query"
$ws.Range("C5").Value = "getChatRoomInList"
$ws.Range("D5").Value = "post, messages, deal"
$ws.Range("E5").Value = "chatId"
$ws.Range("F5").Value = "Chat(+isReviewed)"

# --- Row 6 (mutation/exitChatRoom) is no longer needed: remove it entirely ---
$ws.Rows("6").Delete()

# Match the page setup recorded for the sheet after the edit
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Cursor ends up on E9 after the table shrank
$ws.Range("E9").Select()
